# preparation publication 0.2.0 be6a807bbdadc24333e2c553161780cb6e805524
#
# Changes applied on the "Metadata" worksheet (sheet1):
#   - Version value: 0.1.1 -> 0.2.0
#   - Date value: 2023-10-19T16:17:18+00:00 -> 2023-10-19T17:05:12+00:00
#   - A new "Jurisdiction" / "iso:code:3166:FR" row inserted right after the
#     "Contact" row (pushing Description/Purpose/Copyright/Immutable down by one row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the Version property value.
$ws.Range("B3").Value = "0.2.0"

# Update the publication Date property value.
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# Insert a new row for "Jurisdiction" after the existing "Contact" row (row 10),
# copying the formatting of the row above it so the new row matches the rest
# of the table (borders / wrap text / vertical alignment).
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
